$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Steffan): BB (kg) and TB (cm) corrections, status gizi normal -> obese
$ws.Range("D2").Value = 68
$ws.Range("E2").Value = 150
$ws.Range("K2").Value = "obese"

# Row 3 (Para Gilbert): TB (cm) correction
$ws.Range("E3").Value = 150

# Row 4 (Min Max): TB (cm) corrected from cm to meters
$ws.Range("E4").Value = 1.78

# Row 5 (Vietti): TB (cm) corrected from cm to meters
$ws.Range("E5").Value = 1.84

# Row 6 (Zendaia): TB (cm) corrected from cm to meters
$ws.Range("E6").Value = 1.78

# Row 7 (Zeptr): TB (cm) corrected from cm to meters
$ws.Range("E7").Value = 1.69

# Row 8 (Boiz): TB (cm) corrected from cm to meters, status gizi overweight -> obese
$ws.Range("E8").Value = 1.94
$ws.Range("K8").Value = "obese"

Write-Output "Edits applied"
